# Daily attendance processing - 2025-12-30 01:36:21
# Swap the order of the two comma-separated "recorded by" values in column G
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
